$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=2; B="Bitcoin"; C="https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; D="28.077.87"; E="  +0.14%  "},
    @{Row=3; B="Ethereum"; C="https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; D="1.876.24"; E="  -1.29%  "},
    @{Row=4; B="TetherUSD"; C="https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; D="'1.003"; E="  +0.14%  "},
    @{Row=5; B="BNB"; C="https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D="'313.63"; E="  +0.35%  "},
    @{Row=6; B="USDC"; C="https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; D="'1.002"; E="  +0.10%  "},
    @{Row=7; B="XRP"; C="https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; D="'0.5070"; E="  -0.19%  "},
    @{Row=8; B="Cardano"; C="https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D="'0.3846"; E="  -1.93%  "},
    @{Row=9; B="Dogecoin"; C="https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D="'0.08973"; E="  -2.82%  "},
    @{Row=10; B="Polygon"; C="https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D="'1.123"; E="  -0.88%  "},
    @{Row=11; B="OKB"; C="https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"; D="'41.60"; E="  -0.39%  "},
    @{Row=12; B="Polkadot"; C="https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D="'6.339"; E="  -0.35%  "},
    @{Row=13; B="Solana"; C="https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D="'20.72"; E="  -0.18%  "},
    @{Row=14; B="WrappedEther"; C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D="1.875.42"; E="  -1.26%  "},
    @{Row=15; B="Chainlink"; C="https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D="'7.205"; E="  -1.11%  "},
    @{Row=16; B="BinanceUSD"; C="https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; D="'1.003"; E="  +0.12%  "},
    @{Row=17; B="ShibaInu"; C="https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D="'0.00001107"; E="  -0.81%  "},
    @{Row=18; B="Litecoin"; C="https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D="'91.16"; E="  -1.28%  "},
    @{Row=19; B="TRON"; C="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D="'0.06597"; E="  +0.25%  "},
    @{Row=20; B="Avalanche"; C="https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D="'18.14"; E="  +2.10%  "},
    @{Row=21; B="Dai"; C="https://coinranking.com/coin/MoTuySvg7+dai-dai"; D="'1.002"; E="  +0.16%  "},
    @{Row=22; B="Uniswap"; C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D="'6.110"; E="  -1.67%  "},
    @{Row=23; B="WrappedBTC"; C="https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D="28.098.59"; E="  +0.02%  "},
    @{Row=24; B="Cosmos"; C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D="'11.41"; E="  +0.56%  "},
    @{Row=25; B="Toncoin"; C="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D="'2.274"; E="  -1.94%  "},
    @{Row=26; B="WrappedliquidstakedEther2.0"; C="https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D="2.093.37"; E="  -1.12%  "},
    @{Row=27; B="LidoDAOToken"; C="https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D="'2.536"; E="  -2.24%  "},
    @{Row=28; B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="'20.75"; E="  -0.69%  "},
    @{Row=29; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="'156.91"; E="  -0.30%  "},
    @{Row=30; B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="'126.70"; E="  -0.38%  "},
    @{Row=31; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="'0.1051"; E="  -1.59%  "},
    @{Row=32; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="'1.059"; E="  -2.40%  "},
    @{Row=33; B="Filecoin"; C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D="'5.606"; E="  +0.08%  "},
    @{Row=34; B="HuobiToken"; C="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D="'3.600"; E="  -0.33%  "},
    @{Row=35; B="FraxShare"; C="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D="'9.605"; E="  +0.29%  "},
    @{Row=36; B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="'0.06580"; E="  -1.03%  "},
    @{Row=37; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="'0.02420"; E="  +0.59%  "},
    @{Row=38; B="Algorand"; C="https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; D="'0.2175"; E="  +0.31%  "},
    @{Row=39; B="TrustWalletToken"; C="https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D="'1.272"; E="  +1.13%  "},
    @{Row=40; B="ARBITRUM"; C="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D="'1.207"; E="  -1.41%  "},
    @{Row=41; B="TheSandbox"; C="https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; D="'0.6397"; E="  +0.87%  "},
    @{Row=42; B="Aptos"; C="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D="'11.47"; E="  +0.70%  "},
    @{Row=43; B="InternetComputer(DFINITY)"; C="https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D="'4.917"; E="  -1.03%  "},
    @{Row=44; B="Decentraland"; C="https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"; D="'0.6028"; E="  +0.91%  "},
    @{Row=45; B="EnergySwap"; C="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D="'13.14"; E="  -1.03%  "},
    @{Row=46; B="PancakeSwap"; C="https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D="'3.675"; E="  -0.73%  "},
    @{Row=47; B="WEMIXTOKEN"; C="https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; D="'1.276"; E="  -0.09%  "},
    @{Row=48; B="EOS"; C="https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"; D="'1.239"; E="  +5.23%  "},
    @{Row=49; B="NEARProtocol"; C="https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"; D="'1.995"; E="  -0.60%  "},
    @{Row=50; B="Quant"; C="https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; D="'121.29"; E="  -0.89%  "},
    @{Row=51; B="Aave"; C="https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D="'79.72"; E="  +2.18%  "}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}
